# Generate Report for Handback
# The f9b94de5-... file has been handed back (in sync with en-US) for both
# the zh-cn and de-de locales. Update the Overview sheet's status columns,
# each locale sheet's Status column, and stamp the "Latest Handback DateTime"
# for the row that was just handed back.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row for f9b94de5-...md (row 3) flips from
#     "Ready for handoff" to "Handed back: in sync with en-US" for both
#     locale columns (B = zh-cn, C = de-de).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: Status column (C) for the f9b94de5 row (row 3) updates,
#     and the Latest Handback DateTime column (H) is stamped for both rows.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("H2").Value = "2016-03-15 04:07:38"
$wsZhCn.Range("H3").Value = "2016-03-15 04:07:38"

# --- de-de sheet: same update pattern with its own handback timestamp.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("H2").Value = "2016-03-15 04:07:51"
$wsDeDe.Range("H3").Value = "2016-03-15 04:07:51"
